$d = $word.ActiveDocument

[void]$d.Content.Find.Execute("48+0=48", $true, $false, $false, $false, $false, $true, 1, $false, "22-2=20", 2)
[void]$d.Content.Find.Execute("16+0=16", $true, $false, $false, $false, $false, $true, 1, $false, "78-6=72", 2)
[void]$d.Content.Find.Execute("70-61=9", $true, $false, $false, $false, $false, $true, 1, $false, "5+51=56", 2)
[void]$d.Content.Find.Execute("16+56=72", $true, $false, $false, $false, $false, $true, 1, $false, "74-46=28", 2)
[void]$d.Content.Find.Execute("64-11=53", $true, $false, $false, $false, $false, $true, 1, $false, "56-48=8", 2)
[void]$d.Content.Find.Execute("4+59=63", $true, $false, $false, $false, $false, $true, 1, $false, "39+16=55", 2)
[void]$d.Content.Find.Execute("67-51=16", $true, $false, $false, $false, $false, $true, 1, $false, "36-15=21", 2)
[void]$d.Content.Find.Execute("4+12=16", $true, $false, $false, $false, $false, $true, 1, $false, "88-52=36", 2)
[void]$d.Content.Find.Execute("35+45=80", $true, $false, $false, $false, $false, $true, 1, $false, "26+11=37", 2)
[void]$d.Content.Find.Execute("9-3=6", $true, $false, $false, $false, $false, $true, 1, $false, "37+8=45", 2)
[void]$d.Content.Find.Execute("92-10=82", $true, $false, $false, $false, $false, $true, 1, $false, "33+11=44", 2)
[void]$d.Content.Find.Execute("71-4=67", $true, $false, $false, $false, $false, $true, 1, $false, "23+76=99", 2)
[void]$d.Content.Find.Execute("68-63=5", $true, $false, $false, $false, $false, $true, 1, $false, "22+73=95", 2)
[void]$d.Content.Find.Execute("19+56=75", $true, $false, $false, $false, $false, $true, 1, $false, "33+63=96", 2)
[void]$d.Content.Find.Execute("19+64=83", $true, $false, $false, $false, $false, $true, 1, $false, "93-77=16", 2)
[void]$d.Content.Find.Execute("60+5=65", $true, $false, $false, $false, $false, $true, 1, $false, "3-3=0", 2)
[void]$d.Content.Find.Execute("3-1=2", $true, $false, $false, $false, $false, $true, 1, $false, "70-1=69", 2)
[void]$d.Content.Find.Execute("15+52=67", $true, $false, $false, $false, $false, $true, 1, $false, "24+55=79", 2)
[void]$d.Content.Find.Execute("2+80=82", $true, $false, $false, $false, $false, $true, 1, $false, "6+75=81", 2)
[void]$d.Content.Find.Execute("6+47=53", $true, $false, $false, $false, $false, $true, 1, $false, "40+25=65", 2)
[void]$d.Content.Find.Execute("53+20=73", $true, $false, $false, $false, $false, $true, 1, $false, "11-6=5", 2)
[void]$d.Content.Find.Execute("66-46=20", $true, $false, $false, $false, $false, $true, 1, $false, "0+40=40", 2)
[void]$d.Content.Find.Execute("56+12=68", $true, $false, $false, $false, $false, $true, 1, $false, "36+55=91", 2)
[void]$d.Content.Find.Execute("27+24=51", $true, $false, $false, $false, $false, $true, 1, $false, "14+2=16", 2)
[void]$d.Content.Find.Execute("27+25=52", $true, $false, $false, $false, $false, $true, 1, $false, "9+36=45", 2)
[void]$d.Content.Find.Execute("91-39=52", $true, $false, $false, $false, $false, $true, 1, $false, "89-14=75", 2)
[void]$d.Content.Find.Execute("66-23=43", $true, $false, $false, $false, $false, $true, 1, $false, "85+3=88", 2)
[void]$d.Content.Find.Execute("89-73=16", $true, $false, $false, $false, $false, $true, 1, $false, "6+0=6", 2)
[void]$d.Content.Find.Execute("14-4=10", $true, $false, $false, $false, $false, $true, 1, $false, "76-59=17", 2)
[void]$d.Content.Find.Execute("6+12=18", $true, $false, $false, $false, $false, $true, 1, $false, "3+65=68", 2)
[void]$d.Content.Find.Execute("21+69=90", $true, $false, $false, $false, $false, $true, 1, $false, "59-17=42", 2)
[void]$d.Content.Find.Execute("62+7=69", $true, $false, $false, $false, $false, $true, 1, $false, "84-73=11", 2)
[void]$d.Content.Find.Execute("30-21=9", $true, $false, $false, $false, $false, $true, 1, $false, "83-39=44", 2)
[void]$d.Content.Find.Execute("42-37=5", $true, $false, $false, $false, $false, $true, 1, $false, "80+8=88", 2)
[void]$d.Content.Find.Execute("53-4=49", $true, $false, $false, $false, $false, $true, 1, $false, "12+32=44", 2)
[void]$d.Content.Find.Execute("80-66=14", $true, $false, $false, $false, $false, $true, 1, $false, "45-24=21", 2)
[void]$d.Content.Find.Execute("86-75=11", $true, $false, $false, $false, $false, $true, 1, $false, "73-48=25", 2)
[void]$d.Content.Find.Execute("36+8=44", $true, $false, $false, $false, $false, $true, 1, $false, "47-27=20", 2)
[void]$d.Content.Find.Execute("92-74=18", $true, $false, $false, $false, $false, $true, 1, $false, "82+13=95", 2)
[void]$d.Content.Find.Execute("67+16=83", $true, $false, $false, $false, $false, $true, 1, $false, "82-0=82", 2)
[void]$d.Content.Find.Execute("72+2=74", $true, $false, $false, $false, $false, $true, 1, $false, "59-2=57", 2)
[void]$d.Content.Find.Execute("85-56=29", $true, $false, $false, $false, $false, $true, 1, $false, "11+2=13", 2)
[void]$d.Content.Find.Execute("60+22=82", $true, $false, $false, $false, $false, $true, 1, $false, "6+13=19", 2)
[void]$d.Content.Find.Execute("46-37=9", $true, $false, $false, $false, $false, $true, 1, $false, "57-49=8", 2)
[void]$d.Content.Find.Execute("3+56=59", $true, $false, $false, $false, $false, $true, 1, $false, "75-31=44", 2)
[void]$d.Content.Find.Execute("34+24=58", $true, $false, $false, $false, $false, $true, 1, $false, "11+21=32", 2)
[void]$d.Content.Find.Execute("4+9=13", $true, $false, $false, $false, $false, $true, 1, $false, "79+12=91", 2)
[void]$d.Content.Find.Execute("23-14=9", $true, $false, $false, $false, $false, $true, 1, $false, "83-7=76", 2)
[void]$d.Content.Find.Execute("42-20=22", $true, $false, $false, $false, $false, $true, 1, $false, "97-9=88", 2)
[void]$d.Content.Find.Execute("35+31=66", $true, $false, $false, $false, $false, $true, 1, $false, "82-5=77", 2)
[void]$d.Content.Find.Execute("14+4=18", $true, $false, $false, $false, $false, $true, 1, $false, "80-34=46", 2)
[void]$d.Content.Find.Execute("17+46=63", $true, $false, $false, $false, $false, $true, 1, $false, "51+13=64", 2)
[void]$d.Content.Find.Execute("76-13=63", $true, $false, $false, $false, $false, $true, 1, $false, "30+67=97", 2)
[void]$d.Content.Find.Execute("25-5=20", $true, $false, $false, $false, $false, $true, 1, $false, "4+19=23", 2)
[void]$d.Content.Find.Execute("70+29=99", $true, $false, $false, $false, $false, $true, 1, $false, "86-51=35", 2)
[void]$d.Content.Find.Execute("20+76=96", $true, $false, $false, $false, $false, $true, 1, $false, "29+54=83", 2)
[void]$d.Content.Find.Execute("29-5=24", $true, $false, $false, $false, $false, $true, 1, $false, "2+20=22", 2)
[void]$d.Content.Find.Execute("54-4=50", $true, $false, $false, $false, $false, $true, 1, $false, "14+72=86", 2)
[void]$d.Content.Find.Execute("32+23=55", $true, $false, $false, $false, $false, $true, 1, $false, "55-42=13", 2)
[void]$d.Content.Find.Execute("25+10=35", $true, $false, $false, $false, $false, $true, 1, $false, "41+51=92", 2)
[void]$d.Content.Find.Execute("39+30=69", $true, $false, $false, $false, $false, $true, 1, $false, "45-40=5", 2)
[void]$d.Content.Find.Execute("90-67=23", $true, $false, $false, $false, $false, $true, 1, $false, "41-14=27", 2)
[void]$d.Content.Find.Execute("95-34=61", $true, $false, $false, $false, $false, $true, 1, $false, "94-51=43", 2)
[void]$d.Content.Find.Execute("53+22=75", $true, $false, $false, $false, $false, $true, 1, $false, "44-16=28", 2)
[void]$d.Content.Find.Execute("68+30=98", $true, $false, $false, $false, $false, $true, 1, $false, "45-44=1", 2)
[void]$d.Content.Find.Execute("51+47=98", $true, $false, $false, $false, $false, $true, 1, $false, "2+31=33", 2)
[void]$d.Content.Find.Execute("57+24=81", $true, $false, $false, $false, $false, $true, 1, $false, "74-58=16", 2)
[void]$d.Content.Find.Execute("22+45=67", $true, $false, $false, $false, $false, $true, 1, $false, "19+22=41", 2)
[void]$d.Content.Find.Execute("88-1=87", $true, $false, $false, $false, $false, $true, 1, $false, "47+23=70", 2)
[void]$d.Content.Find.Execute("73+2=75", $true, $false, $false, $false, $false, $true, 1, $false, "51+3=54", 2)
[void]$d.Content.Find.Execute("0+32=32", $true, $false, $false, $false, $false, $true, 1, $false, "38+35=73", 2)
[void]$d.Content.Find.Execute("26-20=6", $true, $false, $false, $false, $false, $true, 1, $false, "63-15=48", 2)
[void]$d.Content.Find.Execute("6-0=6", $true, $false, $false, $false, $false, $true, 1, $false, "52+20=72", 2)
[void]$d.Content.Find.Execute("38-29=9", $true, $false, $false, $false, $false, $true, 1, $false, "41+18=59", 2)
[void]$d.Content.Find.Execute("77-7=70", $true, $false, $false, $false, $false, $true, 1, $false, "65+19=84", 2)
[void]$d.Content.Find.Execute("90-64=26", $true, $false, $false, $false, $false, $true, 1, $false, "27+54=81", 2)
[void]$d.Content.Find.Execute("74-50=24", $true, $false, $false, $false, $false, $true, 1, $false, "51+30=81", 2)
[void]$d.Content.Find.Execute("19+17=36", $true, $false, $false, $false, $false, $true, 1, $false, "88-64=24", 2)
[void]$d.Content.Find.Execute("29+19=48", $true, $false, $false, $false, $false, $true, 1, $false, "59-19=40", 2)
[void]$d.Content.Find.Execute("13+26=39", $true, $false, $false, $false, $false, $true, 1, $false, "47+41=88", 2)
[void]$d.Content.Find.Execute("43+55=98", $true, $false, $false, $false, $false, $true, 1, $false, "30-4=26", 2)
[void]$d.Content.Find.Execute("70-55=15", $true, $false, $false, $false, $false, $true, 1, $false, "24-2=22", 2)
[void]$d.Content.Find.Execute("12-8=4", $true, $false, $false, $false, $false, $true, 1, $false, "97-60=37", 2)
[void]$d.Content.Find.Execute("68-65=3", $true, $false, $false, $false, $false, $true, 1, $false, "89-87=2", 2)
[void]$d.Content.Find.Execute("69-61=8", $true, $false, $false, $false, $false, $true, 1, $false, "97-60=37", 2)
[void]$d.Content.Find.Execute("3+37=40", $true, $false, $false, $false, $false, $true, 1, $false, "17+40=57", 2)
[void]$d.Content.Find.Execute("50+7=57", $true, $false, $false, $false, $false, $true, 1, $false, "72+26=98", 2)
[void]$d.Content.Find.Execute("72+6=78", $true, $false, $false, $false, $false, $true, 1, $false, "37-10=27", 2)
[void]$d.Content.Find.Execute("86+9=95", $true, $false, $false, $false, $false, $true, 1, $false, "45-39=6", 2)
[void]$d.Content.Find.Execute("62+3=65", $true, $false, $false, $false, $false, $true, 1, $false, "46-10=36", 2)
[void]$d.Content.Find.Execute("75-29=46", $true, $false, $false, $false, $false, $true, 1, $false, "29+46=75", 2)
[void]$d.Content.Find.Execute("93-43=50", $true, $false, $false, $false, $false, $true, 1, $false, "24+14=38", 2)
[void]$d.Content.Find.Execute("31+20=51", $true, $false, $false, $false, $false, $true, 1, $false, "8+82=90", 2)
[void]$d.Content.Find.Execute("10+9=19", $true, $false, $false, $false, $false, $true, 1, $false, "54+34=88", 2)
[void]$d.Content.Find.Execute("69+21=90", $true, $false, $false, $false, $false, $true, 1, $false, "58+37=95", 2)
[void]$d.Content.Find.Execute("24+50=74", $true, $false, $false, $false, $false, $true, 1, $false, "64+19=83", 2)
[void]$d.Content.Find.Execute("48+14=62", $true, $false, $false, $false, $false, $true, 1, $false, "4+94=98", 2)
[void]$d.Content.Find.Execute("19+31=50", $true, $false, $false, $false, $false, $true, 1, $false, "51-12=39", 2)
[void]$d.Content.Find.Execute("71+28=99", $true, $false, $false, $false, $false, $true, 1, $false, "36+43=79", 2)
[void]$d.Content.Find.Execute("27-14=13", $true, $false, $false, $false, $false, $true, 1, $false, "58+11=69", 2)
